$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order IDs refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555675267024"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555695848374"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255569586839"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555696484015"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651255569720862"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555674960942.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555675100548.csv"
$ws1.Range("B4").Value = "go_stims-1651255567511055.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555675257034.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-1651255567636856.csv"
$ws2.Range("B3").Value = "OB-16512555679069927.csv"
$ws2.Range("B4").Value = "TB-1651255569563162.csv"
$ws2.Range("B5").Value = "ZB-match_5-16512555675697546.csv"
$ws2.Range("B6").Value = "ZB-match_5-16512555676640444.csv"
$ws2.Range("B7").Value = "TB-1651255568804532.csv"
$ws2.Range("B8").Value = "TB-16512555693977525.csv"
$ws2.Range("B9").Value = "OB-16512555678440838.csv"
$ws2.Range("B10").Value = "OB-16512555677369988.csv"

# --- Sheet 3: RS --- (no cell data changes, only name changed above)

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555696156812.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555695937574.csv"
$ws4.Range("B4").Value = "MM_stims-1651255569631853.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555696166816.csv"
$ws4.Range("B6").Value = "MM_stims-1651255569647399.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555696327174.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555696742113.csv"
$ws5.Range("B3").Value = "SAT_stims-1651255569654583.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555696902192.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555697051325.csv"
